$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.186.72'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.35%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.856.81'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.87%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.34'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.14%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4769'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.83%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2822'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06519'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.864.25'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.35%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07350'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.04%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.41'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.154'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '87.34'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.90%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6466'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.54%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '30.149.01'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.28%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.25'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007622'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.113.90'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.282'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '217.98'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +14.25%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.120'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.319'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.91%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.16'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.77%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.52'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.86%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.912'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.41%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.426'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.52%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.261'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.57%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09122'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.73%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.977'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.08%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05044'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.34%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7447'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.35%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.137'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.99%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.691'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.88%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01825'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.54%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.610'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.09%  '
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9070'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.79%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.050'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.38%  '
$ws.Range("B41").Value = 'Quant'
$ws.Range("C41").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '107.16'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.07%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.914'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.70%  '
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4260'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.54%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.002'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.79%  '
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.443'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.00%  '
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1314'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.07%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.568'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +9.69%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '64.33'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -8.06%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.854'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.47%  '
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.30'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.04%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05705'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.18%  '
